$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("caseType1-valid-values")
$ws.Name = "caseType1-vl"
